# "Diseno de pruebas clase Pelota"
# - Remove the stray "esHoja" test row from the PruebasUnit sheet (it was not
#   an actual method of the Pelota class).
# - Mark the "Diseno" column (E) with "x" for every remaining Pelota test row,
#   since the test design for that class is now complete.
# - Update the active selection to reflect where the author was working.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PruebasUnit")
$ws.Activate()

# Delete the "esHoja" row (row 15) entirely; rows below shift up by one.
$ws.Rows("15:15").Delete()

# Mark "Diseno" (column E) as done for the whole Pelota block (now rows 14-21).
$ws.Range("E14:E21").Value = "x"

# Reflect the author's on-screen selection at save time.
$ws.Range("E27:E33").Select()
